$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data
# (e.g. thousand-separated "29.504.90" or zero-padded "0.07760"). Excel
# auto-coerces single-decimal numeric-looking text into real numbers on
# assignment, which would silently drop trailing zeros / exponent layout,
# so those specific cells are pre-formatted as Text before the assignment.

$ws.Range("D2").Value = "29.504.90"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "1.855.65"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.32"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6941"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3063"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07663"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.52"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "1.855.68"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.138"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6933"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.76"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.296"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "29.524.95"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008279"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "2.107.11"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.34"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.71"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.621"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1482"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.905"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.54"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.22"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.534"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.245"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.139"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05224"
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7773"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.873"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.146"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "1.313.65"
$ws.Range("E38").Value = "  +6.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01868"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9440"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.07"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.802"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.755"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.006.62"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5229"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000122"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.780"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.77"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05951"
$ws.Range("E51").Value = "  +0.67%  "
